$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1705426356589147
$ws.Range("C2").Value = 0.5930232558139535
$ws.Range("J2").Value = 0.007751937984496124
$ws.Range("P2").Value = 0.1279069767441861
$ws.Range("S2").Value = 0.1007751937984496
$ws.Range("C3").Value = 0.0130718954248366
$ws.Range("J3").Value = 0.006535947712418301
$ws.Range("P3").Value = 0.7843137254901961
$ws.Range("S3").Value = 0.196078431372549
$ws.Range("J4").Value = 0.02439024390243903
$ws.Range("P4").Value = 0.7073170731707317
$ws.Range("S4").Value = 0.2682926829268293
$ws.Range("B6").Value = 0.03043478260869565
$ws.Range("D6").Value = 0.01739130434782609
$ws.Range("F6").Value = 0.05217391304347826
$ws.Range("J6").Value = 0.208695652173913
$ws.Range("O6").Value = 0.01739130434782609
$ws.Range("Q6").Value = 0.1739130434782609
$ws.Range("R6").Value = 0.08695652173913043
$ws.Range("S6").Value = 0.4130434782608696
$ws.Range("B7").Value = 0.09574468085106383
$ws.Range("D7").Value = 0.03191489361702127
$ws.Range("F7").Value = 0.05851063829787234
$ws.Range("J7").Value = 0.1276595744680851
$ws.Range("O7").Value = 0.01063829787234043
$ws.Range("Q7").Value = 0.1702127659574468
$ws.Range("R7").Value = 0.0851063829787234
$ws.Range("S7").Value = 0.4202127659574468
$ws.Range("B8").Value = 0.07660455486542443
$ws.Range("D8").Value = 0.02070393374741201
$ws.Range("F8").Value = 0.05383022774327122
$ws.Range("J8").Value = 0.1221532091097308
$ws.Range("O8").Value = 0.02070393374741201
$ws.Range("Q8").Value = 0.2132505175983437
$ws.Range("R8").Value = 0.09109730848861283
$ws.Range("S8").Value = 0.401656314699793
$ws.Range("B9").Value = 0.09417040358744394
$ws.Range("D9").Value = 0.004484304932735426
$ws.Range("F9").Value = 0.1076233183856502
$ws.Range("J9").Value = 0.08968609865470852
$ws.Range("O9").Value = 0.02242152466367713
$ws.Range("Q9").Value = 0.179372197309417
$ws.Range("R9").Value = 0.08071748878923767
$ws.Range("S9").Value = 0.42152466367713
$ws.Range("B10").Value = 0.1012658227848101
$ws.Range("D10").Value = 0.01661392405063291
$ws.Range("E10").Value = 0.002373417721518987
$ws.Range("F10").Value = 0.07041139240506329
$ws.Range("J10").Value = 0.1036392405063291
$ws.Range("O10").Value = 0.01028481012658228
$ws.Range("Q10").Value = 0.2294303797468354
$ws.Range("R10").Value = 0.1004746835443038
$ws.Range("S10").Value = 0.3655063291139241
$ws.Range("G11").Value = 0.1603053435114504
$ws.Range("J11").Value = 0.05725190839694656
$ws.Range("K11").Value = 0.1870229007633588
$ws.Range("L11").Value = 0.583969465648855
$ws.Range("S11").Value = 0.01145038167938931
$ws.Range("G12").Value = 0.75
$ws.Range("J12").Value = 0.16875
$ws.Range("K12").Value = 0.0125
$ws.Range("L12").Value = 0.04375
$ws.Range("S12").Value = 0.025
$ws.Range("G13").Value = 0.6727272727272727
$ws.Range("J13").Value = 0.2181818181818182
$ws.Range("S13").Value = 0.1090909090909091
$ws.Range("F15").Value = 0.0184331797235023
$ws.Range("H15").Value = 0.1751152073732719
$ws.Range("I15").Value = 0.06451612903225806
$ws.Range("J15").Value = 0.3548387096774194
$ws.Range("K15").Value = 0.05990783410138249
$ws.Range("M15").Value = 0.0184331797235023
$ws.Range("O15").Value = 0.07373271889400922
$ws.Range("S15").Value = 0.2350230414746544
$ws.Range("F16").Value = 0.03278688524590164
$ws.Range("H16").Value = 0.1967213114754098
$ws.Range("I16").Value = 0.09289617486338798
$ws.Range("J16").Value = 0.3770491803278688
$ws.Range("K16").Value = 0.1311475409836066
$ws.Range("M16").Value = 0.02185792349726776
$ws.Range("O16").Value = 0.0273224043715847
$ws.Range("S16").Value = 0.1202185792349727
$ws.Range("F17").Value = 0.005882352941176471
$ws.Range("H17").Value = 0.1529411764705882
$ws.Range("I17").Value = 0.09607843137254903
$ws.Range("J17").Value = 0.4529411764705882
$ws.Range("K17").Value = 0.08431372549019608
$ws.Range("M17").Value = 0.02352941176470588
$ws.Range("O17").Value = 0.07647058823529412
$ws.Range("S17").Value = 0.107843137254902
$ws.Range("F18").Value = 0.02608695652173913
$ws.Range("H18").Value = 0.1869565217391304
$ws.Range("I18").Value = 0.07391304347826087
$ws.Range("J18").Value = 0.4173913043478261
$ws.Range("K18").Value = 0.09130434782608696
$ws.Range("M18").Value = 0.01304347826086956
$ws.Range("O18").Value = 0.08695652173913043
$ws.Range("S18").Value = 0.1043478260869565
$ws.Range("F19").Value = 0.01401869158878505
$ws.Range("H19").Value = 0.2227414330218069
$ws.Range("I19").Value = 0.09890965732087227
$ws.Range("J19").Value = 0.3613707165109034
$ws.Range("K19").Value = 0.08489096573208722
$ws.Range("M19").Value = 0.02414330218068536
$ws.Range("N19").Value = 0.000778816199376947
$ws.Range("O19").Value = 0.06619937694704049
$ws.Range("S19").Value = 0.1269470404984424
